$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 55
$ws.Range("D5").Value = 3.1415926

$ws.Range("C8").Select()
